# Fix incorrect data in the "Authors" column (E) for rows 2-4.
# The author-list strings used a run of spaces after each comma as a
# padding marker; the cleaned data bumps that run from 15 spaces to 17
# spaces (old separator ",<15 spaces>" -> new separator ",<17 spaces>").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$spaces15 = ""
for ($i = 0; $i -lt 15; $i++) { $spaces15 = $spaces15 + " " }
$spaces17 = ""
for ($i = 0; $i -lt 17; $i++) { $spaces17 = $spaces17 + " " }

$oldSep = "," + $spaces15
$newSep = "," + $spaces17

foreach ($addr in @("E2", "E3", "E4")) {
    $cell = $ws.Range($addr)
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldSep, $newSep)
}
